# Insert a new weekly record row at row 338 (pushes existing rows 338:430 down to 339:431)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("338:338").Insert()

$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 45093
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = 100112012
$ws.Range("G338").Value = "Espinaca"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 1500
$ws.Range("K338").Value = 450
$ws.Range("L338").Value = 500
$ws.Range("M338").Value = 475
$ws.Range("N338").Value = "$/atado 300 a 500 gramos"
$ws.Range("O338").Value = "Provincia del Elquí"
$ws.Range("P338").Value = 950
$ws.Range("Q338").Value = 0.5
$ws.Range("R338").Value = "Hortaliza"
